# Remove the rows for courses that are no longer offered:
#   - COMPUTO FLEXIBLE (SOFTCOMPUTING)                         (row 9)
#   - PROYECTO DE GESTION DE LA TECNOLOGIA DE INFORMACION      (row 23)
#   - PROYECTO DE SISTEMAS ROBUSTOS, PARALELOS Y DISTRIBUIDOS  (row 24)
# Deleting entire rows shifts the remaining data up and shrinks the
# used range from A1:D38 down to A1:D35, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so row numbers of rows still to be removed
# don't shift while we work.
$ws.Rows.Item(24).EntireRow.Delete()
$ws.Rows.Item(23).EntireRow.Delete()
$ws.Rows.Item(9).EntireRow.Delete()
